# Re-analyzed carbohydrate esters & amides
#
# The "litterChemistry" sheet had one combined category for carbohydrate
# esters ("lipid" spectral range) and one for amides ("amide" spectral
# range). Re-analysis split each back into two original sub-categories
# (carboEster1/carboEster2 and amide1/amide2), each gets its own row of
# partial eta-squared values, inserted around the existing lipid/alkane/
# amide rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("litterChemistry")

# --- 1. Make room for the new rows -----------------------------------
# Before:           After:
#  row2 glycosidicBond     row2 glycosidicBond
#  row3 C_O_stretching     row3 C_O_stretching
#  row4 carboEster         row4 carboEster
#  row5 lipid               row5 carboEster1   (new)
#  row6 alkane               row6 carboEster2   (new)
#  row7 amide                row7 lipid         (was row5)
#                             row8 alkane        (was row6)
#                             row9 amide         (was row7)
#                             row10 amide1       (new)
#                             row11 amide2       (new)

# Insert two new blank rows right before the old "lipid" row (row 5),
# pushing lipid/alkane/amide down to rows 7/8/9.
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()

# Insert two new blank rows right after the (now shifted) "amide" row
# (row 9), for the amide1/amide2 rows.
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(10).Insert()

# --- 2. Fill in category labels + values for every data row ----------
# (rows 2-4, 7-9 keep their original values; only set here for safety /
# clarity since Insert() can disturb formatting.)

function Set-Row {
    param($r, $cat, $b, $c, $d, $f, $g)
    $ws.Cells.Item($r, 1).Value = $cat
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 6).Value = $f
    $ws.Cells.Item($r, 7).Value = $g
}

Set-Row 2  "glycosidicBond" $null     0.69689900000000005 0.083413000000000001 $null 0.067702999999999999
Set-Row 3  "C_O_stretching" $null     0.58411199999999996 $null                $null $null
Set-Row 4  "carboEster"     $null     0.56182100000000001 $null                $null $null
Set-Row 5  "carboEster1"    $null     0.54624799999999996 $null                $null $null
Set-Row 6  "carboEster2"    $null     0.39500000000000002 0.20715700000000001  $null 0.14499400000000001
Set-Row 7  "lipid"          $null     0.78331499999999998 $null                $null $null
Set-Row 8  "alkane"         0.21379300000000001 0.63507400000000003 $null      0.29460399999999998 $null
Set-Row 9  "amide"          $null     $null                0.21581700000000001 $null $null
Set-Row 10 "amide1"         0.15523300000000001 $null       0.197297            $null $null
Set-Row 11 "amide2"         $null     $null                0.21951499999999999 $null $null

# --- 3. Give the 4 newly-inserted rows the same thin-border look as
#        the rest of the data rows (Insert() leaves them unstyled). ---
foreach ($r in 5, 6, 10, 11) {
    for ($c = 1; $c -le 8; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.Borders.LineStyle = 1
        $cell.Borders.Weight = 2
    }
}

# --- 4. Match the saved selection/active cell -------------------------
$ws.Range("D11").Select() | Out-Null

$wb.Save()
